$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 13923
$ws1.Range("F6").Value = 495
$ws1.Range("F7").Value = 1214
$ws1.Range("F8").Value = 1039
$ws1.Range("F9").Value = 13907
$ws1.Range("F10").Value = 14862
$ws1.Range("F12").Value = 8
$ws1.Range("F24").Value = 121
$ws1.Range("F26").Value = 5776
$ws1.Range("F29").Value = 5429
$ws1.Range("F30").Value = 53
$ws1.Range("F31").Value = 54
$ws1.Range("F32").Value = 297

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 13923
$ws4.Range("F7").Value = 495
$ws4.Range("F8").Value = 1214
$ws4.Range("F9").Value = 1039
$ws4.Range("F10").Value = 13907
$ws4.Range("F11").Value = 14862
$ws4.Range("F13").Value = 8
$ws4.Range("F25").Value = 121
$ws4.Range("F27").Value = 5776
$ws4.Range("F30").Value = 5429
$ws4.Range("F31").Value = 53
$ws4.Range("F32").Value = 54
$ws4.Range("F33").Value = 297
